# Apply the OOXML changes described by the diff to the "Overall Results"
# bar chart on slide 7 ("Diagramm 6" graphic frame), plus the frame's
# resize/reposition on the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)            # "Diagramm 6" chart graphic frame
$chart = $sh.Chart

# --- 1. Chart title font color: Text 1, Lighter 35% -> plain Text 1 (black) ---
$title = $chart.ChartTitle
$titleFont = $title.TextFrame.TextRange.Font
$titleFont.Color.RGB = 0            # RGB(0,0,0) packed BGR = 0x000000

# --- 2. Series 1 (bars) fill: Accent 2, Darker 25% -> literal RGB 0076A3 ---
$series = $chart.SeriesCollection(1)
$series.Format.Fill.ForeColor.RGB = 0xA37600   # packed BGR for #0076A3

# --- 3. Value axis: set a fixed maximum scale of 1 ---
$valueAxis = $chart.Axes(2, 1)      # xlValue, xlPrimary
$valueAxis.MaximumScale = 1

# --- 4. Resize / reposition the chart graphic frame on the slide ---
$emuPerPt = 12700
$halfEmuInPt = 0.5 / $emuPerPt
$sh.Left   = 2204382 / $emuPerPt + $halfEmuInPt
$sh.Top    = 1372281 / $emuPerPt + $halfEmuInPt
$sh.Width  = 8684392 / $emuPerPt + $halfEmuInPt
$sh.Height = 4608389 / $emuPerPt + $halfEmuInPt

# --- 5. Rename the shape to match the re-uploaded copy ---
$sh.Name = "Diagramm 7"
